$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30; this shifts existing rows 30-81 down to 31-82,
# carrying all their data (and formats) with them.
$ws.Rows.Item(30).EntireRow.Insert()

# Populate the newly inserted row 30 with the new weekly record. The non-varying
# "template" columns (market/category/etc.) match every other row in the block.
$ws.Cells.Item(30, 1).Value = 10
$ws.Cells.Item(30, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(30, 3).Value = "La Araucanía"
$ws.Cells.Item(30, 4).Value = 44469
$ws.Cells.Item(30, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(30, 5).Value = 9
$ws.Cells.Item(30, 6).Value = 100112012
$ws.Cells.Item(30, 7).Value = "Espinaca"
$ws.Cells.Item(30, 8).Value = "Sin especificar"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 20
$ws.Cells.Item(30, 11).Value = 12000
$ws.Cells.Item(30, 12).Value = 12000
$ws.Cells.Item(30, 13).Value = 12000
$ws.Cells.Item(30, 14).Value = "$/docena de atados"
$ws.Cells.Item(30, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(30, 16).Value = 4000
$ws.Cells.Item(30, 17).Value = 3
$ws.Cells.Item(30, 18).Value = "Hortaliza"
